$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 34, shifting existing rows 34:66 down to 35:67
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new weekly price record
$ws.Cells.Item(34, 1).Value = 5
$ws.Cells.Item(34, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(34, 3).Value = "Maule"
$ws.Cells.Item(34, 4).Value = 44467
$ws.Cells.Item(34, 5).Value = 7
$ws.Cells.Item(34, 6).Value = 100112001
$ws.Cells.Item(34, 7).Value = "Berenjena"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 300
$ws.Cells.Item(34, 11).Value = 7000
$ws.Cells.Item(34, 12).Value = 7000
$ws.Cells.Item(34, 13).Value = 7000
$ws.Cells.Item(34, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(34, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(34, 16).Value = 117
$ws.Cells.Item(34, 17).Value = 60
$ws.Cells.Item(34, 18).Value = "Hortaliza"
